# Foxlocket timings workbook update
# "Working, timer divider = 1024"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item('Sheet1')

# ---------------------------------------------------------------------
# 1. Defined names: add Addr (Sheet1!$B$16) and Multi (Sheet1!$B$17)
# ---------------------------------------------------------------------
$wb.Names.Add('Addr', '=Sheet1!$B$16')
$wb.Names.Add('Multi', '=Sheet1!$B$17')

# ---------------------------------------------------------------------
# 2. Update existing values
# ---------------------------------------------------------------------
$ws.Range('B5').Value = 64
$ws.Range('B7').Value = 200
$ws.Range('B10').Value = 8.4

# ---------------------------------------------------------------------
# 3. Row 11 becomes a computed "Packet duration" row (was the plain
#    "Timer1 divider" / literal-64 row with a stray D11 note)
# ---------------------------------------------------------------------
$ws.Range('D11').ClearContents()

$ws.Range('A11').Value = 'Packet duration'
$ws.Range('A11').Style = 'Calculation'

$ws.Range('B11').Formula = '=B10/B6'
$ws.Range('B11').Style = 'Calculation'

# C11 already holds the shared "tics" string with no explicit style - leave it.

# ---------------------------------------------------------------------
# 4. Remove old rows 12-14 ("Address" / "WhenTransmit" / "WhenTransmit")
#    - their content is superseded by the new block below.
# ---------------------------------------------------------------------
$ws.Range('A12:D14').Clear()

# ---------------------------------------------------------------------
# 5. New shared-string order in the source workbook is: Time, Multiplier,
#    Addr 1..Addr 10 - seed the "Time" text first (row 21/22 block) so the
#    shared-string table indices line up, then fill in the rest of the
#    rows in their natural sheet order.
# ---------------------------------------------------------------------
$ws.Range('A21').Value = 'Time'
$ws.Range('A22').Value = 'Time'

# ---------------------------------------------------------------------
# 6. New "Address" / "Multiplier" / "WhenTransmit" block (rows 16-19)
# ---------------------------------------------------------------------
$ws.Range('A16').Value = 'Address'
$ws.Range('A16').Style = 'Good'
$ws.Range('B16').Value = 1
$ws.Range('B16').Style = 'Good'

$ws.Range('A17').Value = 'Multiplier'
$ws.Range('A17').Style = 'Good'
$ws.Range('B17').Value = 128
$ws.Range('B17').Style = 'Good'

$ws.Range('A18').Value = 'WhenTransmit'
$ws.Range('A18').Style = 'Calculation'
$ws.Range('B18').Formula = '=Addr*B17+B11'
$ws.Range('B18').Style = 'Calculation'
$ws.Range('C18').Value = 'tics'

$ws.Range('A19').Value = 'WhenTransmit'
$ws.Range('B19').Formula = '=B18*B6'
$ws.Range('C19').Value = 'ms'

# ---------------------------------------------------------------------
# 7. Finish the "Time" block (rows 21-22) - text already seeded above
# ---------------------------------------------------------------------
$ws.Range('A21').Style = 'Good'
$ws.Range('B21').Value = 7.4
$ws.Range('B21').Style = 'Good'
$ws.Range('C21').Value = 'ms'

$ws.Range('A22').Style = 'Calculation'
$ws.Range('B22').Formula = '=B21/B6'
$ws.Range('B22').Style = 'Calculation'
$ws.Range('C22').Value = 'tics'

# ---------------------------------------------------------------------
# 8. "Addr 1" .. "Addr 10" multiplication table (rows 25-34)
# ---------------------------------------------------------------------
$addrLabels = @('Addr 1', 'Addr 2', 'Addr 3', 'Addr 4', 'Addr 5', 'Addr 6', 'Addr 7', 'Addr 8', 'Addr 9', 'Addr 10')
for ($i = 0; $i -lt 10; $i++) {
    $r = 25 + $i
    $ws.Cells.Item($r, 1).Value = $addrLabels[$i]
    $ws.Cells.Item($r, 2).Value = $i + 1
    $ws.Cells.Item($r, 3).Formula = '=B' + $r + '*Multi'
}

# ---------------------------------------------------------------------
# 9. Sheet view / selection
# ---------------------------------------------------------------------
$ws.Range('B18').Select()

# ---------------------------------------------------------------------
# 10. Page setup
# ---------------------------------------------------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
